# Update "想去人数" (F column) counts that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3460
$wsExhibit.Range("F4").Value = 137
$wsExhibit.Range("F5").Value = 6998
$wsExhibit.Range("F6").Value = 2515
$wsExhibit.Range("F7").Value = 49
$wsExhibit.Range("F11").Value = 82
$wsExhibit.Range("F14").Value = 581

# Sheet "全部类型" (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3460
$wsAll.Range("F5").Value = 137
$wsAll.Range("F6").Value = 6998
$wsAll.Range("F7").Value = 2515
$wsAll.Range("F8").Value = 49
$wsAll.Range("F12").Value = 82
$wsAll.Range("F15").Value = 581
